$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 5
$ws.Range("B6").Value = 0
